$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.924.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.772.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.44%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.61"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.31%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.789.29"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.80"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.81%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.263.85"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.873.23"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.96%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.787.23"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.10%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.563"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.64%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.177"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.17%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0969"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +15.60%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.37"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +10.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.27%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.13%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "342.74"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.75%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.99"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0610"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.20%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.15%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.177.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.31%  "
